$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (model/header names): H and I swap meaning ---
$ws.Range("H1").Value = "SamplePortion"
$ws.Range("I1").Value = "Result"

# --- Row 2 (data types): H and I now carry a unit annotation ---
$ws.Range("H2").Value = "#float,  unit:mg"
$ws.Range("I2").Value = "#float,  unit:mg"

# --- Row 3 (new): human-readable French descriptions / enum for each column ---
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#PriseEssai"
$ws.Range("I3").Value = "#Resultat"
$ws.Range("J3").Value = "#NuméroLotReactif"
